$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old wide table (A1:M2)
$ws.Range("A1:M2").ClearContents()

# New headers: A1 = "t", B1 = "Chuva"
$ws.Range("A1").Value = "t"
$ws.Range("B1").Value = "Chuva"

# Data: t values 1-12 in column A, Chuva values in column B, rows 2-13
$tValues = 1..12
$chuvaValues = @(233, 199, 137, 84, 100, 101, 104, 84, 144, 148, 134, 168)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tValues[$i]
    $ws.Cells.Item($row, 2).Value = $chuvaValues[$i]
}

# Update selection to match target (activeCell F18, sqref F18)
$ws.Range("F18").Select()
